# Full scale massive run
# Bump several "constant" sheet distribution parameters for a full-scale run,
# and leave the view positioned on the last edited row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("constant")

# --- Update parameter values (distri_param1 / distri_param2 columns) ---
# clinical_progression_rate (row 2): distri_param2 (E2) 5 -> 10
$ws.Range("E2").Value = 10

# progression_rate (row 14): distri_param2 (E14) 2 -> 3
$ws.Range("E14").Value = 3

# raw_transmission_rate (row 17): distri_param1 (D17) 10 -> 1
$ws.Range("D17").Value = 1

# infectiousness_gain_rate (row 21): distri_param2 (E21) 5 -> 10
$ws.Range("E21").Value = 10

# --- Update view/selection state to match where the author left off ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 9
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E27").Select()
